$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'65.863.88"
$ws.Range("E2").Value = "  -2.32%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.429.94"
$ws.Range("E3").Value = "  -1.09%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'581.72"
$ws.Range("E5").Value = "  -2.03%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'173.43"
$ws.Range("E6").Value = "  -4.47%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = "  -3.73%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "'3.433.69"
$ws.Range("E9").Value = "  -1.00%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -7.50%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -2.12%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.410"
$ws.Range("E12").Value = "  -4.66%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "'4.030.04"
$ws.Range("E13").Value = "  -0.91%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.16%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "'29.91"
$ws.Range("E15").Value = "  -6.52%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "'65.966.34"
$ws.Range("E16").Value = "  -2.18%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -4.04%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "'3.429.95"
$ws.Range("E18").Value = "  -1.17%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -5.34%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'13.72"
$ws.Range("E20").Value = "  -2.79%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'365.94"
$ws.Range("E21").Value = "  -7.26%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -3.17%  "

# Row 23 - Dai
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.27%  "

# Row 24 & 25 - LEO / Litecoin swap positions
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'71.81"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "'5.71"
$ws.Range("E25").Value = "  -1.48%  "

# Row 26 - Polygon
$ws.Range("D26").Value = "'0.527"
$ws.Range("E26").Value = "  -2.10%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  -3.63%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  -6.77%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  +0.98%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.00%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "'23.96"
$ws.Range("E31").Value = "  +1.79%  "

# Row 32 - NEARProtocol
$ws.Range("D32").Value = "'5.73"
$ws.Range("E32").Value = "  -6.30%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  -3.75%  "

# Row 34 - USDe
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.05%  "

# Row 35 - Fetch.AI
$ws.Range("D35").Value = "'1.28"
$ws.Range("E35").Value = "  -8.68%  "

# Row 36 - Aptos
$ws.Range("D36").Value = "'7.00"
$ws.Range("E36").Value = "  -4.58%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -2.62%  "

# Row 38 - Monero
$ws.Range("D38").Value = "'159.70"
$ws.Range("E38").Value = "  -0.82%  "

# Row 39 - EnergySwap
$ws.Range("D39").Value = "'28.91"
$ws.Range("E39").Value = "  +10.27%  "

# Row 40 - Mantle
$ws.Range("D40").Value = "'0.877"
$ws.Range("E40").Value = "  -1.38%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -4.66%  "

# Row 42 - dogwifhat
$ws.Range("D42").Value = "'2.56"
$ws.Range("E42").Value = "  -10.01%  "

# Row 43 - Maker
$ws.Range("D43").Value = "'2.711.25"
$ws.Range("E43").Value = "  -1.19%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  -6.19%  "

# Row 45 - RenderToken
$ws.Range("D45").Value = "'6.26"
$ws.Range("E45").Value = "  -7.53%  "

# Row 46 - Hedera
$ws.Range("D46").Value = "'0.0679"
$ws.Range("E46").Value = "  -5.39%  "

# Row 47 - OKB
$ws.Range("D47").Value = "'39.87"
$ws.Range("E47").Value = "  -4.03%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "'24.01"
$ws.Range("E48").Value = "  -8.43%  "

# Row 49 - VeChain
$ws.Range("D49").Value = "'0.0287"
$ws.Range("E49").Value = "  -3.91%  "

# Row 50 - Bittensor
$ws.Range("D50").Value = "'303.22"
$ws.Range("E50").Value = "  -6.97%  "

# Row 51 - SuiNetwork
$ws.Range("E51").Value = "  -3.92%  "
